$d = $word.ActiveDocument

# Phase 1: move each distinct content block to a unique placeholder token
# (breaks the rotation cycle so no search text is ever destroyed before it's read)
$null = $d.Content.Find.Execute("Objetivos: Apresentar os conceitos básicos de Ergonomia e suas aplicações no projeto e operação de Sistemas de Produção e no Desenvolvimento de Produtos.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT0@@", 2)
$null = $d.Content.Find.Execute("Introduce the basic concepts of ergonomics and applications in design and production systems and in the development of products", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT1@@", 2)
$null = $d.Content.Find.Execute("5840917 - Fabricio Maciel Gomes", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT2@@", 2)
$null = $d.Content.Find.Execute("Introdução a ergonomia. Ergonomia Industrial. Ergonomia do Produto. Engenharia de Fatores Humanos", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT3@@", 2)
$null = $d.Content.Find.Execute("Introduction to ergonomics. Industrial Ergonomics. Ergonomics of the product. Human Factors Engineering.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT4@@", 2)
$null = $d.Content.Find.Execute("Conceitos gerais em ergonomia e fatores humanos^lErgonomia Física ^lErgonomia Cognitiva. ^lErgonomia Organizacional. ^lMetodologia de Análise Ergonômica do Trabalho^lErgonomia do Produto. ^lEngenharia de Fatores Humanos", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT5@@", 2)
$null = $d.Content.Find.Execute("Aulas expositivas e práticas.^l", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT6@@", 2)
$null = $d.Content.Find.Execute("Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2^l", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT7@@", 2)
$null = $d.Content.Find.Execute("A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT8@@", 2)
$null = $d.Content.Find.Execute("IIDA, I. Ergonomia: Projeto e Produção. São Paulo: Edgard Blücher, 2001.^lGUÉRIN et al. Compreender o trabalho para transformá-lo. São Paulo: Edgard Blücher, 2001.^lVIDAL, M. C. Ergonomia na empresa, útil, prática e aplicada, 2º ed.. Rio de Janeiro: Editora CVC, 2002.^lCOUTO, H. A. Como implantar a ergonomia na empresa. Belo Horizonte: Ergo Editora, 2002.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT9@@", 2)

# Phase 2: replace each placeholder token with its final destination text
$null = $d.Content.Find.Execute("@@SLOT0@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução a ergonomia. Ergonomia Industrial. Ergonomia do Produto. Engenharia de Fatores Humanos", 2)
$null = $d.Content.Find.Execute("@@SLOT1@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction to ergonomics. Industrial Ergonomics. Ergonomics of the product. Human Factors Engineering.", 2)
$null = $d.Content.Find.Execute("@@SLOT2@@", $true, $false, $false, $false, $false, $true, 1, $false, "Objetivos: Apresentar os conceitos básicos de Ergonomia e suas aplicações no projeto e operação de Sistemas de Produção e no Desenvolvimento de Produtos.", 2)
$null = $d.Content.Find.Execute("@@SLOT3@@", $true, $false, $false, $false, $false, $true, 1, $false, "Conceitos gerais em ergonomia e fatores humanos^lErgonomia Física ^lErgonomia Cognitiva. ^lErgonomia Organizacional. ^lMetodologia de Análise Ergonômica do Trabalho^lErgonomia do Produto. ^lEngenharia de Fatores Humanos", 2)
$null = $d.Content.Find.Execute("@@SLOT4@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introduce the basic concepts of ergonomics and applications in design and production systems and in the development of products", 2)
$null = $d.Content.Find.Execute("@@SLOT5@@", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e práticas.", 2)
$null = $d.Content.Find.Execute("@@SLOT6@@", $true, $false, $false, $false, $false, $true, 1, $false, "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2^l", 2)
$null = $d.Content.Find.Execute("@@SLOT7@@", $true, $false, $false, $false, $false, $true, 1, $false, "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.^l", 2)
$null = $d.Content.Find.Execute("@@SLOT8@@", $true, $false, $false, $false, $false, $true, 1, $false, "IIDA, I. Ergonomia: Projeto e Produção. São Paulo: Edgard Blücher, 2001.^lGUÉRIN et al. Compreender o trabalho para transformá-lo. São Paulo: Edgard Blücher, 2001.^lVIDAL, M. C. Ergonomia na empresa, útil, prática e aplicada, 2º ed.. Rio de Janeiro: Editora CVC, 2002.^lCOUTO, H. A. Como implantar a ergonomia na empresa. Belo Horizonte: Ergo Editora, 2002.", 2)
$null = $d.Content.Find.Execute("@@SLOT9@@", $true, $false, $false, $false, $false, $true, 1, $false, "5840917 - Fabricio Maciel Gomes", 2)
